$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 23489, 1669, 14.07369682444578, 139.9390258789062, 287.9834289550781, 213.9612274169922),
    @(4, 126887, 4473, 28.36731500111782, 171.6033325195312, 1241.945434570312, 706.7743835449219),
    @(5, 14895, 1166, 12.77444253859348, 113.7035751342773, 200.4522705078125, 157.0779228210449),
    @(6, 56773, 2790, 20.34874551971326, 195.8195953369141, 448.7347412109375, 322.2771682739258),
    @(7, 7701, 570, 13.51052631578947, 85.38103485107422, 144.3397369384766, 114.8603858947754),
    @(8, 38520, 2210, 17.42986425339366, 97.80214691162109, 682.30078125, 390.0514640808105),
    @(9, 127734, 3329, 38.37008110543707, 158.8050079345703, 1267.395141601562, 713.1000747680664),
    @(10, 15948, 1108, 14.39350180505415, 77.46805572509766, 297.8793334960938, 187.6736946105957),
    @(11, 10121, 682, 14.84017595307918, 98.10954284667969, 152.2115631103516, 125.1605529785156),
    @(12, 27098, 1885, 14.37559681697613, 94.80092620849609, 455.0233764648438, 274.9121513366699),
    @(14, 12965, 1043, 12.43048897411314, 69.87619018554688, 260.612548828125, 165.2443695068359),
    @(15, 31391, 2280, 13.76798245614035, 118.4514541625977, 413.8274841308594, 266.1394691467285),
    @(16, 88898, 3981, 22.33057020849033, 275.7662353515625, 557.4521484375, 416.6091918945312),
    @(18, 9415, 843, 11.16844602609727, 87.38130187988281, 151.11572265625, 119.2485122680664),
    @(20, 36582, 2823, 12.95855472901169, 186.0671997070312, 385.9921569824219, 286.0296783447266),
    @(23, 25451, 1591, 15.99685732243872, 116.950065612793, 321.0515441894531, 219.000804901123),
    @(24, 30999, 2083, 14.88190110417667, 103.4430923461914, 492.1207275390625, 297.781909942627),
    @(25, 7917, 678, 11.67699115044248, 64.07973480224609, 163.9772338867188, 114.0284843444824),
    @(27, 93394, 5085, 18.36656833824975, 281.6588745117188, 683.6231079101562, 482.6409912109375),
    @(30, 46514, 2524, 18.42868462757528, 86.75746154785156, 835.8033447265625, 461.280403137207),
    @(36, 14721, 1300, 11.32384615384615, 114.6783676147461, 222.7246704101562, 168.7015190124512),
    @(37, 20253, 1421, 14.25263898662913, 68.26117706298828, 441.7120056152344, 254.9865913391113),
    @(38, 20163, 1285, 15.69105058365759, 72.34642791748047, 372.8340454101562, 222.5902366638184),
    @(39, 7811, 728, 10.7293956043956, 79.05519104003906, 142.3079528808594, 110.6815719604492),
    @(42, 27038, 1295, 20.87876447876448, 96.73746490478516, 432.50927734375, 264.6233711242676),
    @(44, 26604, 2174, 12.23735050597976, 155.4207000732422, 335.85791015625, 245.6393051147461),
    @(47, 23789, 1136, 20.94102112676056, 107.7107086181641, 314.014892578125, 210.8628005981445),
    @(54, 15958, 703, 22.69985775248933, 111.2905654907227, 213.9639739990234, 162.627269744873),
    @(55, 8198, 579, 14.15889464594128, 58.04082107543945, 208.79833984375, 133.4195804595947)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowData = $data[$i]
    $ws.Cells.Item($rowNum, 1).Value = $rowData[0]
    $ws.Cells.Item($rowNum, 2).Value = $rowData[1]
    $ws.Cells.Item($rowNum, 3).Value = $rowData[2]
    $ws.Cells.Item($rowNum, 4).Value = $rowData[3]
    $ws.Cells.Item($rowNum, 5).Value = $rowData[4]
    $ws.Cells.Item($rowNum, 6).Value = $rowData[5]
    $ws.Cells.Item($rowNum, 7).Value = $rowData[6]
}